# Switch to new Industry Categories in 2 vars and remove 1 var (#89)
#
# This script updates the "Table B6" mapping column (C) to use the new,
# finer-grained EPS industry categories (instead of the old coarse
# "Other industries" / "Chemicals" / etc. buckets and the now-removed
# "Agriculture" category), and rebuilds the "EoDfIP" summary sheet so it
# lists an elasticity for every one of the new categories (expanding from
# 8 categories / 9 rows to 25 categories / 26 rows).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Table B6: remap the "Matching EPS Industry Category" column (C)
# ---------------------------------------------------------------------
$b6 = $wb.Worksheets.Item("Table B6")

$b6Categories = @{
    3  = "agriculture and forestry 01T03"
    4  = "agriculture and forestry 01T03"
    5  = "oil and gas extraction 06"
    6  = "oil and gas extraction 06"
    7  = "coal mining 05"
    8  = "other mining and quarrying 07T08"
    10 = "energy pipelines and gas processing 352T353"
    11 = "construction 41T43"
    12 = "food beverage and tobacco 10T12"
    13 = "textiles apparel and leather 13T15"
    14 = "textiles apparel and leather 13T15"
    15 = "wood products 16"
    16 = "pulp paper and printing 17T18"
    17 = "pulp paper and printing 17T18"
    18 = "pulp paper and printing 17T18"
    19 = "pulp paper and printing 17T18"
    20 = "refined petroleum and coke 19"
    21 = "refined petroleum and coke 19"
    22 = "chemicals 20"
    23 = "chemicals 20"
    24 = "chemicals 20"
    25 = "chemicals 20"
    26 = "chemicals 20"
    27 = "chemicals 20"
    28 = "rubber and plastic products 22"
    29 = "rubber and plastic products 22"
    30 = "cement and other nonmetallic minerals 239"
    31 = "cement and other nonmetallic minerals 239"
    32 = "cement and other nonmetallic minerals 239"
    33 = "cement and other nonmetallic minerals 239"
    34 = "iron and steel 241"
    35 = "other metals 242"
    36 = "iron and steel 241"
    37 = "other metals 242"
    38 = "other metals 242"
    39 = "metal products except machinery and vehicles 25"
    40 = "other machinery 28"
    41 = "computers and electronics 26"
    42 = "road vehicles 29"
    43 = "nonroad vehicles 30"
    44 = "other manufacturing 31T33"
}

foreach ($r in $b6Categories.Keys) {
    $b6.Cells.Item($r, 3).Value = $b6Categories[$r]
}

# ---------------------------------------------------------------------
# 2. EoDfIP: rebuild the category / elasticity table (A2:B26)
# ---------------------------------------------------------------------
$eo = $wb.Worksheets.Item("EoDfIP")

# Widen column A to fit the new, longer category labels.
$eo.Columns.Item(1).ColumnWidth = 46.6

# Ordered list of the new EPS industry categories that now appear on
# this sheet (was 8 categories / rows 2-9, now 25 categories / rows 2-26).
$eoCategories = @(
    "agriculture and forestry 01T03",
    "coal mining 05",
    "oil and gas extraction 06",
    "other mining and quarrying 07T08",
    "food beverage and tobacco 10T12",
    "textiles apparel and leather 13T15",
    "wood products 16",
    "pulp paper and printing 17T18",
    "refined petroleum and coke 19",
    "chemicals 20",
    "rubber and plastic products 22",
    "glass and glass products 231",
    "cement and other nonmetallic minerals 239",
    "iron and steel 241",
    "other metals 242",
    "metal products except machinery and vehicles 25",
    "computers and electronics 26",
    "appliances and electrical equipment 27",
    "other machinery 28",
    "road vehicles 29",
    "nonroad vehicles 30",
    "other manufacturing 31T33",
    "energy pipelines and gas processing 352T353",
    "water and waste 36T39",
    "construction 41T43"
)

# Categories that have no direct row in Table B6 (new categories split
# out of an existing one with no distinct RFF elasticity) just mirror a
# sibling category's computed elasticity, instead of using the
# SUMPRODUCT/SUMIFS lookup formula. They get the same yellow highlight
# Excel applies to a cell when a formula is typed in referencing another
# cell with an already-formatted number style.
$mirrorFormula = @{
    13 = "=B14"
    19 = "=B18"
}

$lastRow = 1 + $eoCategories.Count
for ($i = 0; $i -lt $eoCategories.Count; $i++) {
    $r = $i + 2
    $cat = $eoCategories[$i]
    $eo.Cells.Item($r, 1).Value = $cat

    if ($mirrorFormula.ContainsKey($r)) {
        $eo.Cells.Item($r, 2).Formula = $mirrorFormula[$r]
        $eo.Cells.Item($r, 2).Interior.Color = 65535
    }
    else {
        $ref = "A$r"
        $formula = '=IFERROR(SUMPRODUCT(--(''Table B6''!$C$3:$C$54=$' + $ref + '),''Table B6''!$B$3:$B$54,''Table A1''!$B$3:$B$54)/SUMIFS(''Table A1''!$B$3:$B$54,''Table B6''!$C$3:$C$54,$' + $ref + '),0)'
        $eo.Cells.Item($r, 2).FormulaArray = $formula
    }
}

# Clear anything left over below the new last row (sheet used to stop at
# row 9; make sure no stray formatting/content survives past row 26).
$oldLastRow = 9
if ($oldLastRow -gt $lastRow) {
    $clearRange = $eo.Range($eo.Cells.Item($lastRow + 1, 1), $eo.Cells.Item($oldLastRow, 2))
    $clearRange.Clear()
}

$wb.Application.Calculate()
